$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36, shifting rows 36:93 down to 37:94.
$ws.Rows("36:36").Insert()

# Populate the new row 36 with this week's data (same style/format as row 37 below it).
$ws.Range("A36").Value = 9
$ws.Range("B36").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 45281
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = 100112010
$ws.Range("G36").Value = "Achicoria"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 70
$ws.Range("K36").Value = 7000
$ws.Range("L36").Value = 8000
$ws.Range("M36").Value = 7500
$ws.Range("N36").Value = "`$/caja 16 unidades"
$ws.Range("O36").Value = "Provincia de Quillota"
$ws.Range("P36").Value = 469
$ws.Range("Q36").Value = 16
$ws.Range("R36").Value = "Hortaliza"
